$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets carry the same exhibition rows and
# both need the refreshed "想去人数" (interest count) figures.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 8030
    $ws.Range("F3").Value = 7668
    $ws.Range("F13").Value = 120
    $ws.Range("F14").Value = 1240
}
